$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Gerber  & Green 2012. FEDAI [Descarga](https://drive.google.com/drive/folders/14HDw0lx7v8cduNtj2XNvvZ5fm_lQ7Z6y?usp=sharing), Barbas 2010 [pdf](https://drive.google.com/file/d/15SqCaheQIA_Eg8Q6CxkkF5Gdt2dPdK1Y/view)"
